# Fix 0 hospitalisation error for India
#
# The "type" column (I) of the data dictionary had stale / incorrect
# per-variable type labels (date, timestamp, integer, nominal, double, ...).
# Every data row should simply read "character" in that column.  A handful
# of rows (8, 9, 10, 11, 21, 22, 23, 30, 44) additionally carried a one-off
# highlight style that needs to be cleared back to the sheet's normal style
# now that the values are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2 through 45 hold data; column I (9) is "type".
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 9).Value = "character"
}

# These rows had a special "flagged" style applied to column I; restore the
# default/normal style now that the underlying value problem is fixed.
$flaggedRows = @(8, 9, 10, 11, 21, 22, 23, 30, 44)
foreach ($r in $flaggedRows) {
    $ws.Cells.Item($r, 9).Style = "Normal"
}

# Restore the view to where the author had scrolled/selected when saving.
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("M42").Select()
